$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: drop R12 from the "Until 2020 included" tallies, and include the two
# new radio-galaxy rows (36, 37) in the reference count.
# ---------------------------------------------------------------------------
$ws.Range("X3").Formula = "=COUNTA(R2,R3,R4,R5,R14,S21,S22,S23,R31,R32,R35,R36,R37)"
$ws.Range("Y3").Formula = "=COUNTA(S22,R31,R32)"

# ---------------------------------------------------------------------------
# Row 5 / Row 6: new "Until 2014 included" sub-sample summary
# (shared string order matters: this must be the first new string written)
# ---------------------------------------------------------------------------
$ws.Range("W5").Value = "Until 2014 included"
$ws.Range("X5").Formula = "=COUNTA(R14,S21,S2,R32,R35,R31)"
$ws.Range("Y5").Formula = "=COUNTA(R31,R32)"

$ws.Range("X6").Formula = "=X2/X5"
$ws.Range("Y6").Formula = "=Y2/Y5"

# ---------------------------------------------------------------------------
# Row 7 / Row 8: add explicit zero S147MHz-limit flags, and a new
# "Until 2015 included" sub-sample summary
# ---------------------------------------------------------------------------
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0

$ws.Range("W8").Value = "Until 2015 included"
$ws.Range("X8").Formula = "=COUNTA(R2:R4,R14,S21,R31,R32,R35)"
$ws.Range("Y8").Formula = "=COUNTA(R31,R32)"

$ws.Range("X9").Formula = "=X2/X8"
$ws.Range("Y9").Formula = "=Y3/Y8"

# ---------------------------------------------------------------------------
# Row 12-14: new "RLQSO until 2020" sub-sample summary (formulas only for
# now; the label text is written further below, after rows 36/37, so the
# shared-string table ends up in the same first-use order as the target).
# ---------------------------------------------------------------------------
$ws.Range("X12").Formula = "=COUNT(B2:B102)-COUNT(B36,B37)"
$ws.Range("Y12").Formula = "=COUNTIF(B2:B112,"">=6"")"

$ws.Range("X13").Formula = "=COUNTA(R2,R3,R4,R5,R14,S21,S22,S23,R31,R32,R35)"
$ws.Range("Y13").Formula = "=COUNTA(S22,R31,R32)"

$ws.Range("X14").Formula = "=X12/X13"
$ws.Range("Y14").Formula = "=Y12/Y13"

# ---------------------------------------------------------------------------
# New rows 36 and 37: two radio-galaxies added to the sample
# ---------------------------------------------------------------------------
$ws.Range("A36").Value = "J0856+0224"
$ws.Range("B36").Value = 5.55
$ws.Range("B36").HorizontalAlignment = -4152
$ws.Range("F36").Value = 899.82
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 86.5
$ws.Range("H36").HorizontalAlignment = -4152
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = -1.18
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("P36").Value = 8.937
$ws.Range("Q36").Value = 2.4
$ws.Range("R36").Value = "Drouart+20"

$ws.Range("A37").Value = "J1530+1049"
$ws.Range("B37").Value = 5.72
$ws.Range("B37").HorizontalAlignment = -4152
$ws.Range("F37").Value = 174.88
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 7.5
$ws.Range("H37").HorizontalAlignment = -4152
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = -1.4
$ws.Range("P37").Value = 15.514
$ws.Range("Q37").Value = 10.825
$ws.Range("R37").Value = "Saxena+18"

# ---------------------------------------------------------------------------
# Row 12 label, written last so this shared string is appended after the
# rows 36/37 strings (matches the target workbook's string table order).
# ---------------------------------------------------------------------------
$ws.Range("W12").Value = "RLQSO until 2020"

# ---------------------------------------------------------------------------
# View state: scroll the window over so column T is at the left edge, and
# move the active selection to Z12 (cosmetic, best-effort).
# ---------------------------------------------------------------------------
$ws.Range("Z12").Select()
